# Apply the refreshed cryptocurrency price / 1h-volume figures to cryptos.xlsx
# (values sourced from coinranking.com by the scheduled GitHub Actions job).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells hold plain text (e.g. "26.751.34", "  -2.66%  "), some of which happen to
# look like ordinary decimal numbers (e.g. "205.91"). Force those to Text number
# format first so Excel stores the literal string instead of silently coercing it
# into a numeric value (which would also eat significant trailing zeros).
function Set-TextValue($cell, $text) {
    if ($text -match '^-?\d+(\.\d+)?$') {
        $cell.NumberFormat = "@"
    }
    $cell.Value = $text
}

Set-TextValue $ws.Range("D2") '26.751.34'
Set-TextValue $ws.Range("E2") '  -2.66%  '
Set-TextValue $ws.Range("D3") '1.557.95'
Set-TextValue $ws.Range("E3") '  -0.67%  '
Set-TextValue $ws.Range("D5") '205.91'
Set-TextValue $ws.Range("E5") '  -1.18%  '
Set-TextValue $ws.Range("D6") '0.485'
Set-TextValue $ws.Range("E6") '  -2.34%  '
Set-TextValue $ws.Range("E7") '  +0.25%  '
Set-TextValue $ws.Range("D8") '21.92'
Set-TextValue $ws.Range("E8") '  -0.24%  '
Set-TextValue $ws.Range("E9") '  -0.72%  '
Set-TextValue $ws.Range("E10") '  -1.45%  '
Set-TextValue $ws.Range("E11") '  -0.26%  '
Set-TextValue $ws.Range("D12") '1.780.27'
Set-TextValue $ws.Range("E12") '  -0.60%  '
Set-TextValue $ws.Range("D13") '1.554.72'
Set-TextValue $ws.Range("E13") '  -0.90%  '
Set-TextValue $ws.Range("D14") '3.73'
Set-TextValue $ws.Range("E14") '  -2.58%  '
Set-TextValue $ws.Range("E15") '  -1.43%  '
Set-TextValue $ws.Range("D16") '26.808.03'
Set-TextValue $ws.Range("E16") '  -2.37%  '
Set-TextValue $ws.Range("D17") '61.50'
Set-TextValue $ws.Range("E17") '  -2.99%  '
Set-TextValue $ws.Range("D18") '214.23'
Set-TextValue $ws.Range("E18") '  +0.02%  '
Set-TextValue $ws.Range("E19") '  +0.44%  '
Set-TextValue $ws.Range("D20") '0.0₃0675'
Set-TextValue $ws.Range("E20") '  -2.06%  '
Set-TextValue $ws.Range("E21") '  +0.23%  '
Set-TextValue $ws.Range("D23") '9.30'
Set-TextValue $ws.Range("E23") '  -2.45%  '
Set-TextValue $ws.Range("E24") '  -1.35%  '
Set-TextValue $ws.Range("D25") '151.75'
Set-TextValue $ws.Range("E25") '  -1.02%  '
Set-TextValue $ws.Range("D26") '6.76'
Set-TextValue $ws.Range("E26") '  +0.38%  '
Set-TextValue $ws.Range("E27") '  -1.32%  '
Set-TextValue $ws.Range("E28") '  +0.21%  '
Set-TextValue $ws.Range("E29") '  -1.61%  '
Set-TextValue $ws.Range("E30") '  -3.99%  '
Set-TextValue $ws.Range("D31") '0.0460'
Set-TextValue $ws.Range("E31") '  -2.25%  '
Set-TextValue $ws.Range("D32") '3.14'
Set-TextValue $ws.Range("D33") '1.385.50'
Set-TextValue $ws.Range("E33") '  +1.63%  '
Set-TextValue $ws.Range("D34") '2.89'
Set-TextValue $ws.Range("E34") '  -1.82%  '
Set-TextValue $ws.Range("E35") '  +0.99%  '
Set-TextValue $ws.Range("E36") '  -0.27%  '
Set-TextValue $ws.Range("E37") '  -4.72%  '
Set-TextValue $ws.Range("E38") '  -3.25%  '
Set-TextValue $ws.Range("E39") '  -2.00%  '
Set-TextValue $ws.Range("E40") '  -3.84%  '
Set-TextValue $ws.Range("E41") '  +0.22%  '
Set-TextValue $ws.Range("D42") '0.999'
Set-TextValue $ws.Range("E42") '  +2.66%  '
Set-TextValue $ws.Range("D43") '5.41'
Set-TextValue $ws.Range("E43") '  +2.41%  '
Set-TextValue $ws.Range("E44") '  +1.47%  '
Set-TextValue $ws.Range("D45") '1.75'
Set-TextValue $ws.Range("E45") '  -2.27%  '
Set-TextValue $ws.Range("D46") '63.00'
Set-TextValue $ws.Range("E46") '  -1.87%  '
Set-TextValue $ws.Range("D47") '1.693.56'
Set-TextValue $ws.Range("E47") '  -0.65%  '
Set-TextValue $ws.Range("D48") '85.47'
Set-TextValue $ws.Range("D49") '0.0₇0971'
Set-TextValue $ws.Range("E49") '  -1.66%  '
Set-TextValue $ws.Range("E50") '  -0.43%  '
Set-TextValue $ws.Range("D51") '0.0945'
Set-TextValue $ws.Range("E51") '  -1.02%  '
